# Apply the trading-results update described by the commit:
# "Trade #13 closed at 2026-02-17 23:53:30 - unknown UNKNOWN +0.000%"
#
# - Summary sheet: refresh aggregate stats (capital, P&L, trade/win counts, win rate)
# - Strategy Status sheet: refresh the MarketMaking strategy row to match
# - All Trades / MarketMaking sheets: append the new closed trade (row 14 / trade #13)

$wb = $excel.ActiveWorkbook

# ---- Summary sheet ----
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1500.32   # Current Capital
$summary.Range("B4").Value = 0.32      # Total P&L $
$summary.Range("B5").Value = 0.49      # Total P&L %
$summary.Range("B6").Value = 13        # Total Trades
$summary.Range("B7").Value = 8         # Winning Trades
$summary.Range("B9").Value = 61.54     # Win Rate %

# ---- Strategy Status sheet (MarketMaking row, row 6) ----
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C6").Value = 100.32
$status.Range("D6").Value = 13
$status.Range("E6").Value = 0.32
$status.Range("F6").Value = 0.32
$status.Range("G6").Value = 61.54

# ---- New trade row (row 14) shared by "All Trades" and "MarketMaking" sheets ----
function Add-TradeRow14($ws) {
    $ws.Cells.Item(14, 1).Value = 13                 # A14 Trade #

    # The Date column holds a literal text string ("2026-02-17"), like the
    # other rows above it - not a date serial. Force text entry so Excel's
    # autodetection doesn't turn it into a date value, then restore the
    # default "Normal" style so no stray number format sticks to the cell.
    $dateCell = $ws.Cells.Item(14, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2026-02-17"                   # B14 Date
    $dateCell.Style = "Normal"

    $ws.Cells.Item(14, 3).Value = "23:53:23"          # C14 Time
    $ws.Cells.Item(14, 4).Value = "MarketMaking"      # D14 Strategy
    $ws.Cells.Item(14, 5).Value = "UP"                # E14 Side
    $ws.Cells.Item(14, 6).Value = 0.91                # F14 Entry Price
    $ws.Cells.Item(14, 7).Value = 0.96                # G14 Exit Price
    $ws.Cells.Item(14, 8).Value = "CLOSED"            # H14 Status
    $ws.Cells.Item(14, 9).Value = 5.4945              # I14 P&L %
    $ws.Cells.Item(14, 10).Value = 0.05               # J14 P&L $
    $ws.Cells.Item(14, 11).Value = 100.32             # K14 Capital After
    $ws.Cells.Item(14, 12).Value = 0                  # L14 Entry Slippage (bps)
    $ws.Cells.Item(14, 13).Value = 0                  # M14 Exit Slippage (bps)
    $ws.Cells.Item(14, 14).Value = 0.6                # N14 Confidence
    $ws.Cells.Item(14, 15).Value = "Normal spread capture: 19600 bps"  # O14 Entry Reason
    $ws.Cells.Item(14, 16).Value = "early_exit"       # P14 Exit Reason
    $ws.Cells.Item(14, 17).Value = 0.13               # Q14 Duration (min)
}

Add-TradeRow14 $wb.Worksheets.Item("All Trades")
Add-TradeRow14 $wb.Worksheets.Item("MarketMaking")
